# resuelto frecuencias verticales y horizontales
#
# Adds a new block (rows 35-45) under the existing "Frecuencias Verticales"
# table that computes / documents the vertical-frequency pivot result
# (T-stamp / T-diff rows) and a small "minutes -> label" lookup table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section title (bold, same style already used by the "Frecuencias
#     Verticales" header at A20) -------------------------------------------
$ws.Range("A35").Value = "RESULTADO FRECUENCIA VERTICAL CON PIVOTE POS 3, LONG 10"
$ws.Range("A35").Font.Bold = $true

# --- "string" row: the horizontal timestamps (as time-of-day fractions) ---
$ws.Range("A37").Value = "string"

$times = @(
    0.28472222222222221,
    0.30555555555555552,
    0.3263888888888889,
    0.34722222222222227,
    0.3888888888888889,
    0.40972222222222227,
    0.4375,
    0.47222222222222227,
    0.49305555555555558,
    0.51388888888888895
)
$col = 2
foreach ($t in $times) {
    $cell = $ws.Cells.Item(37, $col)
    $cell.Value = $t
    $cell.NumberFormat = "h:mm"
    $col = $col + 1
}

# --- "T-stamp" row: same instants expressed as millisecond-of-day values --
$ws.Range("A38").Value = "T-stamp"

$stamps = @(35400000, 37200000, 39000000, 40800000, 44400000, 46200000, 48600000, 51600000, 53400000, 55200000)
$col = 2
foreach ($s in $stamps) {
    $ws.Cells.Item(38, $col).Value = $s
    $col = $col + 1
}

# --- "T-diff" row: successive differences (last column repeats total) -----
$ws.Range("A39").Value = "T-diff"

$diffs = @(1800000, 1800000, 1800000, 3600000, 1800000, 2400000, 3000000, 1800000, 1800000, 55200000)
$col = 2
foreach ($d in $diffs) {
    $ws.Cells.Item(39, $col).Value = $d
    $col = $col + 1
}

# --- small lookup table: minutes label -> TS (ms) value --------------------
# Order matters: it reproduces the order these labels were first typed in
# the authored workbook (and therefore their position in the shared
# strings table).
$ws.Range("A41").Value = "30 min"
$ws.Range("B41").Value = 1800000

$ws.Range("A40").Value = "Min => TS"
$ws.Range("A40").Font.Bold = $true

$ws.Range("A42").Value = "60 min"
$ws.Range("B42").Value = 3600000

$ws.Range("A43").Value = "40 min"
$ws.Range("B43").Value = 2400000

$ws.Range("A44").Value = "50 min"
$ws.Range("B44").Value = 3000000

$ws.Range("A45").Value = "12h20min"
$ws.Range("B45").Value = 55200000

# --- explanatory note (added last, ends up last in shared strings) --------
$ws.Range("A36").Value = "calcula el anterior y posterior, en el caso del final, no deberiamos tomar dicho valor"

# --- viewport: scroll so row 28 is at the top and select I35, matching
#     what the author had in view/selected when they saved -----------------
$null = $ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("I35").Select()
